$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the selection / scroll position on sheet "436" (sheet4)
#    without leaving it as the active tab.
# ------------------------------------------------------------------
$ws436 = $wb.Worksheets.Item("436")
$ws436.Activate()
$ws436.Range("F1").Select()

# ------------------------------------------------------------------
# 2) Update the view state on sheet "Arsenal" (sheet6): scroll so
#    D15 is the top-left visible cell and select M2:M37. This sheet
#    should also end up NOT being the active tab.
# ------------------------------------------------------------------
$wsArsenal = $wb.Worksheets.Item("Arsenal")
$wsArsenal.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 4
$wsArsenal.Range("M2:M37").Select()

# ------------------------------------------------------------------
# 3) Add the new "ars_hist" worksheet after "Arsenal" (last sheet).
#    Adding it last makes it the active sheet/tab, matching the
#    activeTab bump from 5 -> 6.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHist = $wb.Worksheets.Add($null, $lastSheet)
$wsHist.Name = "ars_hist"

# Column B is wide enough to show "Total market value"
$wsHist.Columns.Item(2).ColumnWidth = 17.67

# --- Header row -----------------------------------------------------
# Enter the season labels bottom-up (oldest first as typed into the
# sheet) then the remaining headers, matching the original authoring
# order so shared-string indices line up.
$wsHist.Range("A8").Value = "2015/2016"
$wsHist.Range("A7").Value = "2014/2015"
$wsHist.Range("A6").Value = "2013/2014"
$wsHist.Range("A5").Value = "2012/2013"
$wsHist.Range("A4").Value = "2011/2012"
$wsHist.Range("A3").Value = "2010/2011"
$wsHist.Range("A2").Value = "2009/2010"
$wsHist.Range("A1").Value = "Season"
$wsHist.Range("A9").Value = "2016/2017"
$wsHist.Range("C1").Value = "Matches"
$wsHist.Range("D1").Value = "Wins"
$wsHist.Range("E1").Value = "Draws"
$wsHist.Range("F1").Value = "Losses"
$wsHist.Range("G1").Value = "Points"
$wsHist.Range("I1").Value = ",  "
$wsHist.Range("B1").Value = "Total market value"
$wsHist.Range("H1").Formula = "=A1&I1&B1&I1&C1&I1&D1&I1&E1&I1&F1&I1&G1"

# --- Data rows --------------------------------------------------------
# Row 2: 2009/2010
$wsHist.Range("B2").Value = 281
$wsHist.Range("C2").Value = 38
$wsHist.Range("D2").Value = 23
$wsHist.Range("E2").Value = 6
$wsHist.Range("F2").Value = 9
$wsHist.Range("G2").Formula = "=D2*3+E2*1"

# Row 3: 2010/2011
$wsHist.Range("B3").Value = 304.375
$wsHist.Range("C3").Value = 38
$wsHist.Range("D3").Value = 19
$wsHist.Range("E3").Value = 11
$wsHist.Range("F3").Value = 8
$wsHist.Range("G3:G9").Formula = "=D3*3+E3*1"

# Row 4: 2011/2012
$wsHist.Range("B4").Value = 312.925
$wsHist.Range("C4").Value = 38
$wsHist.Range("D4").Value = 21
$wsHist.Range("E4").Value = 7
$wsHist.Range("F4").Value = 10

# Row 5: 2012/2013
$wsHist.Range("B5").Value = 293.75
$wsHist.Range("C5").Value = 38
$wsHist.Range("D5").Value = 21
$wsHist.Range("E5").Value = 10
$wsHist.Range("F5").Value = 7

# Row 6: 2013/2014
$wsHist.Range("B6").Value = 341.2
$wsHist.Range("C6").Value = 38
$wsHist.Range("D6").Value = 24
$wsHist.Range("E6").Value = 7
$wsHist.Range("F6").Value = 7

# Row 7: 2014/2015
$wsHist.Range("B7").Value = 404.35
$wsHist.Range("C7").Value = 38
$wsHist.Range("D7").Value = 22
$wsHist.Range("E7").Value = 9
$wsHist.Range("F7").Value = 7

# Row 8: 2015/2016
$wsHist.Range("B8").Value = 408.6
$wsHist.Range("C8").Value = 38
$wsHist.Range("D8").Value = 20
$wsHist.Range("E8").Value = 11
$wsHist.Range("F8").Value = 7

# Row 9: 2016/2017
$wsHist.Range("B9").Value = 517
$wsHist.Range("C9").Value = 38
$wsHist.Range("D9").Value = 23
$wsHist.Range("E9").Value = 6
$wsHist.Range("F9").Value = 9

# Final selection on the new sheet
$wsHist.Range("A8").Select()
